# Apply refreshed odds values to Jogos_da_Semana_FlashScore_2025-02-24.xlsx
# (data refresh commit: "Atualizando o arquivo XLSX")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.2
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 6.5
$ws.Range("AC2").Value = 11
# Row 3
$ws.Range("Z3").Value = 8.5
$ws.Range("AC3").Value = 21
# Row 5
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("Y5").Value = 5
$ws.Range("AK5").Value = 23
$ws.Range("AO5").Value = 51
$ws.Range("AP5").Value = 2.03
$ws.Range("AQ5").Value = 1.83
$ws.Range("AR5").Value = 4.4
$ws.Range("AS5").Value = 1.22
# Row 9
$ws.Range("J9").Value = 2.15
$ws.Range("L9").Value = 5.3
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.72
$ws.Range("U9").Value = 1.4
$ws.Range("V9").Value = 2.52
$ws.Range("Y9").Value = 6.1
$ws.Range("AB9").Value = 11
$ws.Range("AD9").Value = 32
$ws.Range("AG9").Value = 19.5
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 12.5
# Row 10
$ws.Range("G10").Value = 2.4
$ws.Range("I10").Value = 3.1
$ws.Range("J10").Value = 3.2
$ws.Range("L10").Value = 3.75
$ws.Range("Y10").Value = 7
$ws.Range("AJ10").Value = 8.5
$ws.Range("AK10").Value = 15
$ws.Range("AL10").Value = 12
$ws.Range("AM10").Value = 34
# Row 11
$ws.Range("J11").Value = 1.91
# Row 12
$ws.Range("J12").Value = 1.95
# Row 15
$ws.Range("G15").Value = 4.2
$ws.Range("H15").Value = 3.35
$ws.Range("I15").Value = 1.87
$ws.Range("J15").Value = 4.5
$ws.Range("K15").Value = 2.07
$ws.Range("L15").Value = 2.47
$ws.Range("U15").Value = 1.45
$ws.Range("V15").Value = 2.6
$ws.Range("W15").Value = 1.98
$ws.Range("X15").Value = 1.75
$ws.Range("Y15").Value = 10
$ws.Range("Z15").Value = 24
$ws.Range("AA15").Value = 15
$ws.Range("AB15").Value = 75
$ws.Range("AC15").Value = 50
$ws.Range("AD15").Value = 60
$ws.Range("AF15").Value = 6.8
$ws.Range("AG15").Value = 18.5
$ws.Range("AH15").Value = 110
$ws.Range("AK15").Value = 8.5
$ws.Range("AL15").Value = 9.25
$ws.Range("AM15").Value = 16.5
$ws.Range("AN15").Value = 18
# Row 16
$ws.Range("G16").Value = 2.2
$ws.Range("I16").Value = 2.9
$ws.Range("J16").Value = 2.88
$ws.Range("W16").Value = 1.57
$ws.Range("X16").Value = 2.25
$ws.Range("AD16").Value = 21
$ws.Range("AG16").Value = 12
$ws.Range("AI16").Value = 126
# Row 17
$ws.Range("G17").Value = 1.7
$ws.Range("I17").Value = 4.75
$ws.Range("J17").Value = 2.38
$ws.Range("L17").Value = 5.5
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 8.5
$ws.Range("Y17").Value = 6
$ws.Range("AA17").Value = 9
$ws.Range("AD17").Value = 34
$ws.Range("AE17").Value = 8.5
$ws.Range("AH17").Value = 67
$ws.Range("AK17").Value = 23
$ws.Range("AL17").Value = 17
# Row 18
$ws.Range("G18").Value = 3.7
$ws.Range("I18").Value = 1.9
$ws.Range("L18").Value = 2.6
$ws.Range("Q18").Value = 2.05
$ws.Range("R18").Value = 1.8
$ws.Range("S18").Value = 3.75
$ws.Range("T18").Value = 1.29
$ws.Range("AK18").Value = 9
$ws.Range("AM18").Value = 17
# Row 19
$ws.Range("G19").Value = 1.44
$ws.Range("H19").Value = 4.25
$ws.Range("I19").Value = 6.2
$ws.Range("J19").Value = 1.93
$ws.Range("K19").Value = 2.32
$ws.Range("L19").Value = 5.9
$ws.Range("O19").Value = 1.22
$ws.Range("P19").Value = 3.45
$ws.Range("Q19").Value = 1.65
$ws.Range("R19").Value = 1.98
$ws.Range("S19").Value = 2.55
$ws.Range("T19").Value = 1.39
$ws.Range("W19").Value = 1.83
$ws.Range("X19").Value = 1.78
$ws.Range("Y19").Value = 7.1
$ws.Range("Z19").Value = 6.9
$ws.Range("AA19").Value = 8.25
$ws.Range("AB19").Value = 9.75
$ws.Range("AC19").Value = 11.5
$ws.Range("AD19").Value = 26
$ws.Range("AE19").Value = 12.5
$ws.Range("AF19").Value = 8.5
$ws.Range("AG19").Value = 18.5
$ws.Range("AH19").Value = 90
$ws.Range("AI19").Value = 700
$ws.Range("AJ19").Value = 17
$ws.Range("AK19").Value = 40
$ws.Range("AL19").Value = 20
$ws.Range("AM19").Value = 120
$ws.Range("AN19").Value = 65
$ws.Range("AO19").Value = 65
# Row 20
$ws.Range("G20").Value = 2.47
$ws.Range("J20").Value = 2.95
$ws.Range("O20").Value = 1.2
$ws.Range("P20").Value = 3.6
$ws.Range("Q20").Value = 1.62
$ws.Range("R20").Value = 2.02
$ws.Range("S20").Value = 2.45
$ws.Range("T20").Value = 1.42
$ws.Range("W20").Value = 1.53
$ws.Range("X20").Value = 2.2
$ws.Range("Y20").Value = 10.75
$ws.Range("Z20").Value = 14
$ws.Range("AC20").Value = 18
$ws.Range("AD20").Value = 23
$ws.Range("AF20").Value = 7
$ws.Range("AI20").Value = 250
# Row 22
$ws.Range("H22").Value = 10.5
$ws.Range("I22").Value = 22
$ws.Range("K22").Value = 4.55
$ws.Range("L22").Value = 14.5
$ws.Range("S22").Value = 1.18
$ws.Range("T22").Value = 4.6
$ws.Range("Y22").Value = 28
$ws.Range("Z22").Value = 13
$ws.Range("AA22").Value = 18.5
$ws.Range("AB22").Value = 8.75
$ws.Range("AD22").Value = 32
$ws.Range("AE22").Value = 65
$ws.Range("AF22").Value = 40
$ws.Range("AG22").Value = 50
$ws.Range("AH22").Value = 110
$ws.Range("AI22").Value = 450
$ws.Range("AJ22").Value = 200
$ws.Range("AK22").Value = 500
$ws.Range("AL22").Value = 120
$ws.Range("AN22").Value = 450
$ws.Range("AO22").Value = 200
# Row 24
$ws.Range("G24").Value = 3.4
$ws.Range("J24").Value = 3.75
$ws.Range("Y24").Value = 17
$ws.Range("AD24").Value = 23
$ws.Range("AE24").Value = 21
$ws.Range("AJ24").Value = 12
# Row 25
$ws.Range("G25").Value = 3.5
$ws.Range("H25").Value = 4.5
$ws.Range("J25").Value = 4
$ws.Range("L25").Value = 2.25
$ws.Range("N25").Value = 23
$ws.Range("Q25").Value = 1.36
$ws.Range("R25").Value = 3.1
$ws.Range("Y25").Value = 21
$ws.Range("Z25").Value = 26
$ws.Range("AA25").Value = 13
$ws.Range("AC25").Value = 23
$ws.Range("AD25").Value = 23
$ws.Range("AE25").Value = 26
$ws.Range("AF25").Value = 10
$ws.Range("AG25").Value = 12
$ws.Range("AH25").Value = 29
# Row 26
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 3.4
$ws.Range("I26").Value = 3.75
$ws.Range("J26").Value = 2.75
$ws.Range("K26").Value = 2.05
$ws.Range("L26").Value = 4.33
$ws.Range("Z26").Value = 9
$ws.Range("AB26").Value = 17
$ws.Range("AE26").Value = 8.5
$ws.Range("AJ26").Value = 9.5
$ws.Range("AK26").Value = 19
$ws.Range("AL26").Value = 13
# Row 28
$ws.Range("Q28").Value = 1.75
$ws.Range("R28").Value = 2.05
# Row 29
$ws.Range("S29").Value = 4
$ws.Range("T29").Value = 1.22
# Row 30
$ws.Range("N30").Value = 8.5
$ws.Range("Q30").Value = 2.1
$ws.Range("R30").Value = 1.7
$ws.Range("S30").Value = 3.75
$ws.Range("T30").Value = 1.25
$ws.Range("W30").Value = 2.05
$ws.Range("X30").Value = 1.7
# Row 31
$ws.Range("G31").Value = 2.55
$ws.Range("I31").Value = 2.6
$ws.Range("J31").Value = 3.1
$ws.Range("L31").Value = 3.1
$ws.Range("M31").Value = 1.02
$ws.Range("N31").Value = 12
$ws.Range("O31").Value = 1.22
$ws.Range("P31").Value = 4
$ws.Range("S31").Value = 2.75
$ws.Range("T31").Value = 1.4
$ws.Range("W31").Value = 1.62
$ws.Range("X31").Value = 2.2
$ws.Range("AB31").Value = 26
$ws.Range("AC31").Value = 21
$ws.Range("AK31").Value = 13
$ws.Range("AL31").Value = 10
$ws.Range("AO31").Value = 26
# Row 32
$ws.Range("G32").Value = 1.75
$ws.Range("I32").Value = 4.5
$ws.Range("J32").Value = 2.38
$ws.Range("L32").Value = 4.5
$ws.Range("M32").Value = 1.03
$ws.Range("N32").Value = 10
$ws.Range("O32").Value = 1.29
$ws.Range("P32").Value = 3.5
$ws.Range("Q32").Value = 1.9
$ws.Range("R32").Value = 1.9
$ws.Range("AK32").Value = 23
$ws.Range("AL32").Value = 15
$ws.Range("AM32").Value = 41
$ws.Range("AN32").Value = 34
# Row 34
$ws.Range("L34").Value = 3.25
$ws.Range("M34").Value = 1.03
$ws.Range("N34").Value = 17
$ws.Range("O34").Value = 1.14
$ws.Range("P34").Value = 5.5
$ws.Range("Q34").Value = 1.48
$ws.Range("R34").Value = 2.6
$ws.Range("AB34").Value = 21
$ws.Range("AE34").Value = 19
$ws.Range("AH34").Value = 29
$ws.Range("AJ34").Value = 15
$ws.Range("AK34").Value = 19
$ws.Range("AO34").Value = 23
